$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.985.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5155"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.19%  "

$ws.Range("E8").Value = "  +2.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07186"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07637"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.892.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.234"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008483"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.022.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.050"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.140.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.382"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("E25").Value = "  +10.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.716"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.904"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.780"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09192"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.228"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7651"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.970"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.261"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.595"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01992"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.063"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.612"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1501"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4807"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.594"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("E50").Value = "  +2.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "
